$wb = $excel.ActiveWorkbook

# Scheduled market-data refresh: update computed Leve profit figures
# (current average prices / profits) across several rows on the ALC,
# ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets. A few rows whose HQ
# (or NQ) branch no longer applies have those cells cleared outright
# rather than zeroed, matching how Excel leaves unused cells blank.

$ws = $wb.Worksheets.Item("ALC")
# Row 5: Met a Sticky End
$ws.Range("H5").Value = 564.8333
$ws.Range("I5").Value = 618
$ws.Range("K5").Value = 618
$ws.Range("M5").Value = -503

# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 1365.174
$ws.Range("I28").Value = 1188.8334
$ws.Range("K28").Value = 1188.8334
$ws.Range("M28").Value = -703.8334

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 6200
$ws.Range("I76").Value = 3666.6667
$ws.Range("K76").Value = 3666.6667
$ws.Range("M76").Value = -3351.6667

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 6200
$ws.Range("I79").Value = 3666.6667
$ws.Range("K79").Value = 3666.6667
$ws.Range("M79").Value = -2574.6667

# Row 98: The Dotted Line
$ws.Range("H98").Value = 3580.6428
$ws.Range("I98").Value = 3602.6365
$ws.Range("K98").Value = 3602.6365
$ws.Range("M98").Value = -2104.6365

# Row 106: Making Your Mark
$ws.Range("H106").Value = 3365.25
$ws.Range("I106").Value = 2727.5
$ws.Range("J106").Value = 4003
$ws.Range("K106").Value = 2727.5
$ws.Range("L106").Value = 4003
$ws.Range("M106").Value = -2096.5
$ws.Range("N106").Value = -5265

# Row 113: Amaro Kart
$ws.Range("H113").Value = 11896.9375
$ws.Range("I113").Value = 8980
$ws.Range("J113").Value = 14165.667
$ws.Range("K113").Value = 8980
$ws.Range("L113").Value = 14165.667
$ws.Range("M113").Value = -5726
$ws.Range("N113").Value = -20673.667

# Row 122: Wishful Inking
$ws.Range("H122").Value = 3580.6428
$ws.Range("I122").Value = 3602.6365
$ws.Range("K122").Value = 10807.9095
$ws.Range("M122").Value = -8357.9095

# Row 131: Mindful Study
$ws.Range("H131").Value = 7994.273
$ws.Range("I131").Value = 1288.5
$ws.Range("J131").Value = 75052
$ws.Range("K131").Value = 3865.5
$ws.Range("L131").Value = 225156
$ws.Range("M131").Value = 1174.5
$ws.Range("N131").Value = -235236

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 1771.9259
$ws.Range("I132").Value = 1336.9166
$ws.Range("K132").Value = 4010.7498
$ws.Range("M132").Value = -1480.7498

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 7229.514
$ws.Range("I32").Value = 3690.1553
$ws.Range("J32").Value = 24336.416
$ws.Range("K32").Value = 3690.1553
$ws.Range("L32").Value = 24336.416
$ws.Range("M32").Value = -3403.1553
$ws.Range("N32").Value = -24910.416

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 65499.668
$ws.Range("I45").Value = 65499.668
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 65499.668
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 2552.25
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2552.25
$ws.Range("K122").Value = 0
$ws.Range("N122").Value = -12556.75
$ws.Range("M122").ClearContents()

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2422.9143
$ws.Range("I132").Value = 2164.0715
$ws.Range("K132").Value = 6492.2145
$ws.Range("M132").Value = -3962.2145

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2343.7173
$ws.Range("I134").Value = 2358.6743
$ws.Range("J134").Value = 2129.3333
$ws.Range("K134").Value = 7076.0229
$ws.Range("L134").Value = 6387.999899999999
$ws.Range("M134").Value = -4541.0229
$ws.Range("N134").Value = -11457.9999

$ws = $wb.Worksheets.Item("CRP")
# Row 9: Shields for the Serpents
$ws.Range("H9").Value = 349666.66
$ws.Range("J9").Value = 349666.66
$ws.Range("L9").Value = 349666.66
$ws.Range("N9").Value = -350002.66

# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 216.18182
$ws.Range("I22").Value = 216.18182
$ws.Range("K22").Value = 216.18182
$ws.Range("M22").Value = 133.81818

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2864.2144
$ws.Range("I132").Value = 2851.4583
$ws.Range("J132").Value = 2940.75
$ws.Range("K132").Value = 8554.374899999999
$ws.Range("L132").Value = 8822.25
$ws.Range("M132").Value = -6024.374899999999
$ws.Range("N132").Value = -13882.25

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 14368.65
$ws.Range("I134").Value = 4742.9443
$ws.Range("K134").Value = 14228.8329
$ws.Range("M134").Value = -11693.8329

$ws = $wb.Worksheets.Item("CUL")
# Row 132: More Mezcal
$ws.Range("H132").Value = 1531
$ws.Range("I132").Value = 1456.9445
$ws.Range("J132").Value = 2197.5
$ws.Range("K132").Value = 13112.5005
$ws.Range("L132").Value = 19777.5
$ws.Range("M132").Value = -10582.5005
$ws.Range("N132").Value = -24837.5

# Row 139: Najoothie
$ws.Range("H139").Value = 1657.591
$ws.Range("I139").Value = 1362.3889
$ws.Range("K139").Value = 4087.1667
$ws.Range("M139").Value = 1052.8333

$ws = $wb.Worksheets.Item("GSM")
# Row 52: It's My Business to Know Things
$ws.Range("H52").Value = 57999
$ws.Range("J52").Value = 57999
$ws.Range("L52").Value = 57999
$ws.Range("N52").Value = -58517

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1168.5555
$ws.Range("J97").Value = 556.75
$ws.Range("L97").Value = 556.75
$ws.Range("N97").Value = -1548.75

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 45456976
$ws.Range("I102").Value = 2234.842
$ws.Range("K102").Value = 2234.842
$ws.Range("M102").Value = -612.8420000000001

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 12411.454
$ws.Range("I126").Value = 14419.889
$ws.Range("J126").Value = 3373.5
$ws.Range("K126").Value = 43259.667
$ws.Range("L126").Value = 10120.5
$ws.Range("M126").Value = -40789.667
$ws.Range("N126").Value = -15060.5

# Row 132: On Board for Lar
$ws.Range("H132").Value = 2977.4827
$ws.Range("I132").Value = 2474.8845
$ws.Range("K132").Value = 7424.6535
$ws.Range("M132").Value = -4894.6535

# Row 133: Pendulums of Our Own
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 3031.8235
$ws.Range("I22").Value = 1637.8
$ws.Range("K22").Value = 1637.8
$ws.Range("M22").Value = -1342.8

# Row 27: Fire and Hide
$ws.Range("H27").Value = 3031.8235
$ws.Range("I27").Value = 1637.8
$ws.Range("K27").Value = 1637.8
$ws.Range("M27").Value = -1530.8

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 998.3333
$ws.Range("I46").Value = 998.3333
$ws.Range("K46").Value = 998.3333
$ws.Range("M46").Value = -810.3333

# Row 122: Hell on Leather
$ws.Range("H122").Value = 5832.7144
$ws.Range("I122").Value = 5908.5713
$ws.Range("J122").Value = 5756.857
$ws.Range("K122").Value = 17725.7139
$ws.Range("L122").Value = 17270.571
$ws.Range("M122").Value = -15275.7139
$ws.Range("N122").Value = -22170.571

$ws = $wb.Worksheets.Item("WVR")
# Row 58: Seeing It Through to the End
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 2485.6072
$ws.Range("J122").Value = 3057
$ws.Range("L122").Value = 9171
$ws.Range("N122").Value = -14071

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 5495.5557
$ws.Range("I126").Value = 4993.5
$ws.Range("K126").Value = 14980.5
$ws.Range("M126").Value = -12510.5
